# PALI_rNPV_Model.xlsx — "Updated Risk-Adjusted Development Costs"
#
# UC Phase 2 committed spend increases from -$15M to -$25M (larger Phase 2
# given the planned large registrational study), and the G&A-to-
# commercialization line increases from -$25M to -$42M (updated annual
# burn schedule). Footnotes for both lines are refreshed to reflect the
# new assumptions. All downstream rNPV / fair-value formulas recalc
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rNPV Model")

# --- UC Phase 2 ($M, committed): -15 -> -25, new footnote ---------------
$ws.Range("B67").Value = -25
$c67 = $ws.Range("C67")
$c67.Value = "Oral small molecule, CRO-run, large registrational study planned N=195"
$c67.Font.Name = "Arial"
$c67.Font.Size = 9
$c67.Font.Italic = $true
$c67.Font.Color = 6710886

# --- G&A to commercialization ($M): -25 -> -42, new footnote ------------
$ws.Range("B71").Value = -42
$c71 = $ws.Range("C71")
$c71.Value = "~$6M/yr 2026-2027 (Phase 2) + ~10M/yr 2028-2030 (Phase 3 & NDA)"
$c71.Font.Name = "Arial"
$c71.Font.Size = 9
$c71.Font.Italic = $true
$c71.Font.Color = 6710886

# Leave the final selection on the cell that was actually edited, as the
# live workbook would after this editing session.
$c71.Select() | Out-Null

$excel.Calculate() | Out-Null
